$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts for Term 1 programs
$ws.Range("C2").Value = 30
$ws.Range("C3").Value = 30
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 10

# Remove the last row (Term 3, BK) which is no longer needed
$ws.Rows(25).Delete()

# Match the reviewer's on-screen view: zoomed in and focused near the data
$excel.ActiveWindow.Zoom = 130
$ws.Range("H10").Select()
